$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feature names and importance values per the diff
$ws.Range("A2").Value = "RSI"
$ws.Range("B2").Value = 0.4294032378541743

$ws.Range("A3").Value = "MACD"
$ws.Range("B3").Value = 0.2889822478476525

$ws.Range("A4").Value = "Signal_line"
$ws.Range("B4").Value = 0.1072619769831947

$ws.Range("A5").Value = "VIX_short"
$ws.Range("B5").Value = 0.03805860443213046

$ws.Range("A6").Value = "close_short"
$ws.Range("B6").Value = 0.03329402888711147

$ws.Range("A7").Value = "VIX"
$ws.Range("B7").Value = 0.03056893575784117

$ws.Range("A8").Value = "close_long"
$ws.Range("B8").Value = 0.02964260696362641

$ws.Range("A9").Value = "VIX_long"
$ws.Range("B9").Value = 0.02530395479255034

$ws.Range("A10").Value = "fedrate"
$ws.Range("B10").Value = 0.01748440648171872
